$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Copy the date-column format from the last existing row (A489) so new rows match
$ws.Cells.Item(489, 1).Copy()
$ws.Range("A490:A517").PasteSpecial(-4122)

# Append the new daily RRPONTSYD observations
$ws.Cells.Item(490, 1).Value = 45187
$ws.Cells.Item(490, 2).Value = 1452.942
$ws.Cells.Item(491, 1).Value = 45188
$ws.Cells.Item(491, 2).Value = 1453.324
$ws.Cells.Item(492, 1).Value = 45189
$ws.Cells.Item(492, 2).Value = 1486.984
$ws.Cells.Item(493, 1).Value = 45190
$ws.Cells.Item(493, 2).Value = 1454.115
$ws.Cells.Item(494, 1).Value = 45191
$ws.Cells.Item(494, 2).Value = 1427.575
$ws.Cells.Item(495, 1).Value = 45194
$ws.Cells.Item(495, 2).Value = 1437.31
$ws.Cells.Item(496, 1).Value = 45195
$ws.Cells.Item(496, 2).Value = 1438.301
$ws.Cells.Item(497, 1).Value = 45196
$ws.Cells.Item(497, 2).Value = 1442.805
$ws.Cells.Item(498, 1).Value = 45197
$ws.Cells.Item(498, 2).Value = 1453.366
$ws.Cells.Item(499, 1).Value = 45198
$ws.Cells.Item(499, 2).Value = 1557.569
$ws.Cells.Item(500, 1).Value = 45201
$ws.Cells.Item(500, 2).Value = 1365.739
$ws.Cells.Item(501, 1).Value = 45202
$ws.Cells.Item(501, 2).Value = 1348.465
$ws.Cells.Item(502, 1).Value = 45203
$ws.Cells.Item(502, 2).Value = 1342.031
$ws.Cells.Item(503, 1).Value = 45204
$ws.Cells.Item(503, 2).Value = 1265.132
$ws.Cells.Item(504, 1).Value = 45205
$ws.Cells.Item(504, 2).Value = 1283.461
$ws.Cells.Item(505, 1).Value = 45209
$ws.Cells.Item(505, 2).Value = 1222.44
$ws.Cells.Item(506, 1).Value = 45210
$ws.Cells.Item(506, 2).Value = 1239.382
$ws.Cells.Item(507, 1).Value = 45211
$ws.Cells.Item(507, 2).Value = 1157.319
$ws.Cells.Item(508, 1).Value = 45212
$ws.Cells.Item(508, 2).Value = 1151.818
$ws.Cells.Item(509, 1).Value = 45215
$ws.Cells.Item(509, 2).Value = 1108.819
$ws.Cells.Item(510, 1).Value = 45216
$ws.Cells.Item(510, 2).Value = 1082.502
$ws.Cells.Item(511, 1).Value = 45217
$ws.Cells.Item(511, 2).Value = 1150.781
$ws.Cells.Item(512, 1).Value = 45218
$ws.Cells.Item(512, 2).Value = 1114.179
$ws.Cells.Item(513, 1).Value = 45219
$ws.Cells.Item(513, 2).Value = 1138.756
$ws.Cells.Item(514, 1).Value = 45222
$ws.Cells.Item(514, 2).Value = 1157.976
$ws.Cells.Item(515, 1).Value = 45223
$ws.Cells.Item(515, 2).Value = 1097.875
$ws.Cells.Item(516, 1).Value = 45224
$ws.Cells.Item(516, 2).Value = 1100.617
$ws.Cells.Item(517, 1).Value = 45225
$ws.Cells.Item(517, 2).Value = 1089.85

# Update the SeriesInfo sheet metadata to reflect the new data pull
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

$wsInfo.Range("B3").NumberFormat = "@"
$wsInfo.Range("B3").Value = "2023-10-27"
$wsInfo.Range("B3").ClearFormats()

$wsInfo.Range("B4").NumberFormat = "@"
$wsInfo.Range("B4").Value = "2023-10-27"
$wsInfo.Range("B4").ClearFormats()

$wsInfo.Range("B7").NumberFormat = "@"
$wsInfo.Range("B7").Value = "2023-10-26"
$wsInfo.Range("B7").ClearFormats()

$wsInfo.Range("B14").Value = "2023-10-26 13:01:02-05"
$wsInfo.Range("B15").Value = 92

